$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force cells with numeric-looking replacement text to remain plain text,
# matching the original workbook where these columns are stored as strings.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '71.443.70'
$ws.Range("E2").Value = '  -0.18%  '
$ws.Range("D3").Value = '3.843.47'
$ws.Range("E3").Value = '  +0.75%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '707.56'
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("D6").Value = '173.29'
$ws.Range("E6").Value = '  -0.94%  '
$ws.Range("D7").Value = '3.842.93'
$ws.Range("E7").Value = '  +0.80%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -0.51%  '
$ws.Range("D10").Value = '0.163'
$ws.Range("E10").Value = '  -0.29%  '
$ws.Range("E11").Value = '  +0.51%  '
$ws.Range("E12").Value = '  -0.25%  '
$ws.Range("E13").Value = '  -1.39%  '
$ws.Range("D14").Value = '37.08'
$ws.Range("E14").Value = '  +1.72%  '
$ws.Range("D15").Value = '4.493.83'
$ws.Range("E15").Value = '  +0.79%  '
$ws.Range("D16").Value = '3.866.52'
$ws.Range("E16").Value = '  +0.78%  '
$ws.Range("D17").Value = '71.468.49'
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("D18").Value = '7.28'
$ws.Range("E18").Value = '  +0.75%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = '0.115'
$ws.Range("E19").Value = '  +0.39%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '17.49'
$ws.Range("E20").Value = '  -1.59%  '
$ws.Range("D21").Value = '498.84'
$ws.Range("E21").Value = '  +3.11%  '
$ws.Range("D22").Value = '10.73'
$ws.Range("E22").Value = '  -1.90%  '
$ws.Range("D23").Value = '0.734'
$ws.Range("E23").Value = '  +2.47%  '
$ws.Range("D24").Value = '85.56'
$ws.Range("E24").Value = '  +1.11%  '
$ws.Range("E25").Value = '  +1.65%  '
$ws.Range("D26").Value = '10.72'
$ws.Range("E26").Value = '  +1.39%  '
$ws.Range("E27").Value = '  -0.89%  '
$ws.Range("D28").Value = '3.996.79'
$ws.Range("E28").Value = '  +0.77%  '
$ws.Range("E29").Value = '  -2.39%  '
$ws.Range("D31").Value = '3.12'
$ws.Range("E31").Value = '  -1.00%  '
$ws.Range("D32").Value = '7.52'
$ws.Range("E32").Value = '  -1.47%  '
$ws.Range("E33").Value = '  -2.60%  '
$ws.Range("D34").Value = '29.49'
$ws.Range("E34").Value = '  -0.56%  '
$ws.Range("E35").Value = '  -5.20%  '
$ws.Range("E36").Value = '  -0.44%  '
$ws.Range("D37").Value = '3.808.77'
$ws.Range("E37").Value = '  +1.15%  '
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  -0.11%  '
$ws.Range("E39").Value = '  +0.28%  '
$ws.Range("D40").Value = '2.36'
$ws.Range("E40").Value = '  +3.31%  '
$ws.Range("E41").Value = '  +5.06%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").Value = '6.03'
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = '3.38'
$ws.Range("E43").Value = '  -2.46%  '
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("D46").Value = '0.000319'
$ws.Range("E46").Value = '  +1.23%  '
$ws.Range("D47").Value = '164.00'
$ws.Range("E47").Value = '  -0.55%  '
$ws.Range("D48").Value = '432.32'
$ws.Range("E48").Value = '  +3.00%  '
$ws.Range("D49").Value = '48.98'
$ws.Range("E49").Value = '  +0.52%  '
$ws.Range("D50").Value = '8.75'
$ws.Range("E50").Value = '  +1.26%  '
$ws.Range("E51").Value = '  -0.57%  '
